$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.394
$ws.Range("C2").Value = 41.996
$ws.Range("D2").Value = 57.39

$ws.Range("B3").Value = 15.252
$ws.Range("C3").Value = 9.057
$ws.Range("D3").Value = 24.309

$ws.Range("B4").Value = 9.481999999999999
$ws.Range("C4").Value = 3.429
$ws.Range("D4").Value = 12.911

$ws.Range("B5").Value = 0.26
$ws.Range("C5").Value = 0.024
$ws.Range("D5").Value = 0.284

$ws.Range("B6").Value = 0.307
$ws.Range("C6").Value = 4.8
$ws.Range("D6").Value = 5.107

$ws.Range("B7").Value = 40.695
$ws.Range("C7").Value = 59.306
$ws.Range("D7").Value = 100.001
